# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted at row 25 (pushing the
# existing rows 25-60 down to 26-61, preserving all of their data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25; Excel shifts rows 25:60 down to 26:61.
$ws.Rows("25").Insert()

# Populate the new row 25 with the new weekly record.
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 44868
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100112031
$ws.Range("G25").Value = "Poroto verde"
$ws.Range("H25").Value = "Magnum"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 27000
$ws.Range("L25").Value = 28000
$ws.Range("M25").Value = 27500
$ws.Range("N25").Value = "`$/malla 25 kilos"
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 1100
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"
